$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting the
# existing Late/Date/Outstanding columns one position to the right.
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet (it picks up tabSelected,
# "Transactions" loses it) and move the selection to R7.
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
